$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.511.01'
$ws.Range('E2').Value = '  +12.72%  '
$ws.Range('D3').Value = '1.831.75'
$ws.Range('E3').Value = '  +9.36%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'230.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.93%  '
$ws.Range('D6').Value = "'0.550"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = "'31.63"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.74%  '
$ws.Range('D9').Value = "'47.17"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.63%  '
$ws.Range('E10').Value = '  +7.31%  '
$ws.Range('D11').Value = "'0.0674"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.18%  '
$ws.Range('D12').Value = "'0.0933"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('D13').Value = '2.093.39'
$ws.Range('E13').Value = '  +9.37%  '
$ws.Range('D14').Value = '1.850.96'
$ws.Range('E14').Value = '  +10.53%  '
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('D16').Value = '34.487.68'
$ws.Range('E16').Value = '  +12.62%  '
$ws.Range('D17').Value = "'10.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = "'4.26"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.09%  '
$ws.Range('D19').Value = "'69.84"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').Value = "'259.90"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.65%  '
$ws.Range('D21').Value = '0.0₃0753'
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = "'10.57"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.58%  '
$ws.Range('D24').Value = "'4.36"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D26').Value = "'158.12"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').Value = "'16.75"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.38%  '
$ws.Range('D28').Value = "'7.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.55%  '
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').Value = "'3.90"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +12.46%  '
$ws.Range('D32').Value = "'0.0520"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.85%  '
$ws.Range('E34').Value = '  +8.72%  '
$ws.Range('D35').Value = '1.552.92'
$ws.Range('E35').Value = '  +4.78%  '
$ws.Range('E36').Value = '  +1.89%  '
$ws.Range('E37').Value = '  +5.58%  '
$ws.Range('D38').Value = "'1.31"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +217.32%  '
$ws.Range('D39').Value = "'0.636"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.84%  '
$ws.Range('E40').Value = '  +6.50%  '
$ws.Range('D41').Value = "'85.07"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('E42').Value = '  +5.37%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = "'2.34"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'0.917"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.05%  '
$ws.Range('D45').Value = "'2.15"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.87%  '
$ws.Range('E46').Value = '  +5.36%  '
$ws.Range('D47').Value = "'1.08"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.89%  '
$ws.Range('D48').Value = '1.983.29'
$ws.Range('E48').Value = '  +9.62%  '
$ws.Range('D49').Value = "'12.45"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +28.20%  '
$ws.Range('D50').Value = "'5.82"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.68%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = "'1.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
